{"js": "// 1) Remove the stray leading space in the run \" \u0433. \" (right after the\n//    \"\u00ab${day}\u00bb ${month} ${year}\" merge field) so it reads \"\u0433. \" instead.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The date-line paragraph (\"\u0433. \u0412\u043b\u0430\u0434\u0438\u0432\u043e\u0441\u0442\u043e\u043a ... \u043e\u0442 \u00ab${day}\u00bb ${month} ${year} \u0433. \")\n// is the 4th paragraph (index 3) in the document.\nconst dateParagraph = paragraphs.items[3];\nconst dateRange = dateParagraph.getRange();\n\n// Scope the search to this paragraph only, since \" \u0433. \" (with surrounding\n// spaces) also appears elsewhere in the document and must stay untouched.\nconst matches = dateRange.search(\" \u0433. \", { matchCase: true, matchWholeWord: false });\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length > 0) {\n  matches.items[0].insertText(\"\u0433. \", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Tiny column-width rebalance on the last row of the 3rd table (the\n//    signature block): 2373 -> 2372 dxa and 2431 -> 2432 dxa (net 0 change).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst signatureTable = tables.items[2];\nconst rows = signatureTable.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst lastRow = rows.items[rows.items.length - 1];\nconst cells = lastRow.cells;\ncells.load(\"items\");\nawait context.sync();\n\n// Columns are 0-indexed: 0=2628, 1=2124, 2=236, 3=2373(->2372), 4=2431(->2432)\ncells.items[3].columnWidth = 2372 / 20; // dxa -> points\ncells.items[4].columnWidth = 2432 / 20; // dxa -> points\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1) Remove the stray leading space in the run \" \u0433. \" (right after the\n#    \"\u00ab${day}\u00bb ${month} ${year}\" merge field) so it reads \"\u0433. \" instead.\n#    That text only occurs once in the whole document as its own run,\n#    inside the 4th paragraph (\"\u0433. \u0412\u043b\u0430\u0434\u0438\u0432\u043e\u0441\u0442\u043e\u043a ... \u043e\u0442 \u00ab${day}\u00bb ${month}\n#    ${year} \u0433. \"), so scope the Find/Replace to that paragraph's range\n#    to avoid touching the unrelated \" \u0433. \" substrings elsewhere in the\n#    document (e.g. inside the \"\u0433. \u0412\u043b\u0430\u0434\u0438\u0432\u043e\u0441\u0442\u043e\u043a\u0430\" address text).\n$dateParagraph = $d.Paragraphs.Item(4)\n$find = $dateParagraph.Range.Find\n$find.ClearFormatting()\n$find.Text = \" \u0433. \"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"\u0433. \"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 1)\n\n# ---------------------------------------------------------------------\n# 2) Tiny column-width rebalance on the last row of the 3rd table (the\n#    signature block): 2373 -> 2372 dxa and 2431 -> 2432 dxa (net 0\n#    change, 20 dxa = 1 point).\n$signatureTable = $d.Tables.Item(3)\n$lastRowIndex = $signatureTable.Rows.Count\n$lastRow = $signatureTable.Rows.Item($lastRowIndex)\n$lastRow.Cells.Item(4).Width = 2372 / 20\n$lastRow.Cells.Item(5).Width = 2432 / 20\n"}
